$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for the "MuSCs" sending cluster (old rows 8-10)
$ws.Rows("8:10").Delete()

# Row 2
$ws.Range("B2").Value = "Mmp9"
$ws.Range("C2").Value = "Cd44"
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.001904333333333333
$ws.Range("H2").Value = 0.005713
$ws.Range("I2").Value = 0.01334809965397277
$ws.Range("J2").Value = 0.01334809965397277
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.487621999999999
$ws.Range("N2").Value = 22.462866
$ws.Range("O2").Value = 0.1384395179233961
$ws.Range("P2").Value = 0.1384395179233961
$ws.Range("Q2").Value = 0.014258928162
$ws.Range("R2").Value = 0.128330353458
$ws.Range("S2").Value = 0.00184790448128944
$ws.Range("T2").Value = 0.00184790448128944

# Row 3
$ws.Range("B3").Value = "Mmp9"
$ws.Range("C3").Value = "Cd44"
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 0.001904333333333333
$ws.Range("H3").Value = 0.005713
$ws.Range("I3").Value = 0.01334809965397277
$ws.Range("J3").Value = 0.01334809965397277
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 31.999428
$ws.Range("N3").Value = 95.998284
$ws.Range("O3").Value = 0.5916411627275552
$ws.Range("P3").Value = 0.5916411627275552
$ws.Range("Q3").Value = 0.060937577388
$ws.Range("R3").Value = 0.548438196492
$ws.Range("S3").Value = 0.007897285199479727
$ws.Range("T3").Value = 0.007897285199479727

# Row 4
$ws.Range("B4").Value = "Mmp9"
$ws.Range("C4").Value = "Cd44"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.001904333333333333
$ws.Range("H4").Value = 0.005713
$ws.Range("I4").Value = 0.01334809965397277
$ws.Range("J4").Value = 0.01334809965397277
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.59882166666667
$ws.Range("N4").Value = 43.796465
$ws.Range("O4").Value = 0.2699193193490487
$ws.Range("P4").Value = 0.2699193193490487
$ws.Range("Q4").Value = 0.02780102272722222
$ws.Range("R4").Value = 0.250209204545
$ws.Range("S4").Value = 0.003602909973203603
$ws.Range("T4").Value = 0.003602909973203603

# Row 5
$ws.Range("B5").Value = "Mmp9"
$ws.Range("C5").Value = "Cd44"
$ws.Range("D5").Value = "ECs"
$ws.Range("G5").Value = 0.1407626666666667
$ws.Range("H5").Value = 0.422288
$ws.Range("I5").Value = 0.9866519003460271
$ws.Range("J5").Value = 0.9866519003460271
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.487621999999999
$ws.Range("N5").Value = 22.462866
$ws.Range("O5").Value = 0.1384395179233961
$ws.Range("P5").Value = 0.1384395179233961
$ws.Range("Q5").Value = 1.053977639712
$ws.Range("R5").Value = 9.485798757407998
$ws.Range("S5").Value = 0.1365916134421066
$ws.Range("T5").Value = 0.1365916134421066

# Row 6
$ws.Range("B6").Value = "Mmp9"
$ws.Range("C6").Value = "Cd44"
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 0.1407626666666667
$ws.Range("H6").Value = 0.422288
$ws.Range("I6").Value = 0.9866519003460271
$ws.Range("J6").Value = 0.9866519003460271
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 31.999428
$ws.Range("N6").Value = 95.998284
$ws.Range("O6").Value = 0.5916411627275552
$ws.Range("P6").Value = 0.5916411627275552
$ws.Range("Q6").Value = 4.504324817088
$ws.Range("R6").Value = 40.538923353792
$ws.Range("S6").Value = 0.5837438775280754
$ws.Range("T6").Value = 0.5837438775280754

# Row 7
$ws.Range("B7").Value = "Mmp9"
$ws.Range("C7").Value = "Cd44"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 0.1407626666666667
$ws.Range("H7").Value = 0.422288
$ws.Range("I7").Value = 0.9866519003460271
$ws.Range("J7").Value = 0.9866519003460271
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 14.59882166666667
$ws.Range("N7").Value = 43.796465
$ws.Range("O7").Value = 0.2699193193490487
$ws.Range("P7").Value = 0.2699193193490487
$ws.Range("Q7").Value = 2.054969067991111
$ws.Range("R7").Value = 18.49472161192
$ws.Range("S7").Value = 0.2663164093758451
$ws.Range("T7").Value = 0.2663164093758451

